$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.799.07"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.048.13"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "227.82"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "59.70"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "0.0835"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "2.350.13"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "14.42"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").Value = "21.46"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "5.50"
$ws.Range("E15").Value = "  +6.30%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "2.051.90"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "37.791.76"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "69.50"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "221.92"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").Value = "168.65"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "18.79"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "1.30"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +8.27%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "6.48"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  +4.21%  "
$ws.Range("D38").Value = "3.49"
$ws.Range("E38").Value = "  +7.65%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "18.37"
$ws.Range("E40").Value = "  +9.00%  "
$ws.Range("D41").Value = "1.526.59"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "97.66"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "4.17"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.10"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "2.238.78"
$ws.Range("E51").Value = "  +0.80%  "
